# Applies the MONTANA_2024 cleaning-script fix:
#  - renames header columns to clean machine-friendly names
#  - normalizes "de"/"del"/"el"/"la"/"los" -> "De"/"Del"/"El"/"La"/"Los" in
#    a handful of place names (title-casing the connector words)
#  - renames "TOTAL" -> "Total" on the grand-total row
#  - removes the trailing metadata rows (191-195), shrinking the used range
#    down to A1:D189

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Municipality / state name capitalization fixes ---
$ws.Range("B16").Value  = "Comitán De Domínguez"
$ws.Range("B21").Value  = "Ocozocoautla De Espinosa"
$ws.Range("A37").Value  = "Ciudad De México"
$ws.Range("B46").Value  = "San Juan Del Río"
$ws.Range("A48").Value  = "Estado De México"
$ws.Range("B48").Value  = "Acambay De Ruíz Castañeda"
$ws.Range("B49").Value  = "Naucalpan De Juárez"
$ws.Range("B55").Value  = "Apaseo El Alto"
$ws.Range("B63").Value  = "San Luis De La Paz"
$ws.Range("B64").Value  = "San Miguel De Allende"
$ws.Range("B65").Value  = "Santa Cruz De Juventino Rosas"
$ws.Range("B66").Value  = "Silao De La Victoria"
$ws.Range("B69").Value  = "Acapulco De Juárez"
$ws.Range("B70").Value  = "Coyuca De Catalán"
$ws.Range("B77").Value  = "Técpan De Galeana"
$ws.Range("B78").Value  = "Zihuatanejo De Azueta"
$ws.Range("B83").Value  = "Autlán De Navarro"
$ws.Range("B87").Value  = "Lagos De Moreno"
$ws.Range("B88").Value  = "San Juan De Los Lagos"
$ws.Range("B89").Value  = "San Miguel El Alto"
$ws.Range("B92").Value  = "Tepatitlán De Morelos"
$ws.Range("B93").Value  = "Unión De Tula"
$ws.Range("A97").Value  = "Michoacán De Ocampo"
$ws.Range("B135").Value = "Amealco De Bonfil"
$ws.Range("B141").Value = "San Juan Del Río"
$ws.Range("B149").Value = "Villa De Ramos"
$ws.Range("A167").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B176").Value = "Poza Rica De Hidalgo"
$ws.Range("B186").Value = "Villa De Cos"

# --- Grand total label ---
$ws.Range("A189").Value = "Total"

# --- Remove trailing metadata/footnote rows (191-195) ---
$ws.Range("A191:A195").EntireRow.Delete()
